$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, exactly as captured in the cryptos.xlsx update diff.
# Every cell in this sheet is stored as text (t="inlineStr" in the source).
# Plain-decimal values (e.g. "2.80") would otherwise be auto-parsed by Excel
# into a number (dropping the trailing zero / losing exact text), so those are
# written with a leading apostrophe to force text, then the cell's style is
# reset from the resulting "quote prefix" style back to a never-touched,
# default-styled donor cell (D26) so no stray style/number-format survives.
$changes = @(
    @('D2', '43.029.00'),
    @('E2', '  -1.93%  '),
    @('D3', '2.252.36'),
    @('E3', '  -2.09%  '),
    @('E4', '  +0.47%  '),
    @('D5', '110.53'),
    @('E5', '  +0.78%  '),
    @('D6', '262.82'),
    @('E6', '  -3.30%  '),
    @('D7', '0.614'),
    @('E7', '  -0.96%  '),
    @('E8', '  +0.33%  '),
    @('D9', '0.595'),
    @('E9', '  -3.65%  '),
    @('D10', '47.11'),
    @('E10', '  -0.14%  '),
    @('D11', '0.0917'),
    @('E11', '  -2.14%  '),
    @('D12', '8.62'),
    @('E12', '  +2.75%  '),
    @('D13', '0.106'),
    @('E13', '  -0.74%  '),
    @('D14', '15.30'),
    @('E14', '  -2.82%  '),
    @('D15', '2.595.23'),
    @('E15', '  -1.71%  '),
    @('D16', '0.845'),
    @('E16', '  -1.54%  '),
    @('D17', '2.259.06'),
    @('E17', '  -1.45%  '),
    @('D18', '42.974.85'),
    @('E18', '  -1.89%  '),
    @('D19', '0.0000106'),
    @('E19', '  -3.88%  '),
    @('D20', '6.86'),
    @('E20', '  +8.71%  '),
    @('D21', '70.54'),
    @('E21', '  -2.36%  '),
    @('D22', '2.38'),
    @('E22', '  -4.55%  '),
    @('B23', 'BitcoinCash'),
    @('C23', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'),
    @('D23', '228.64'),
    @('E23', '  -2.21%  '),
    @('B24', 'InternetComputer(DFINITY)'),
    @('C24', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'),
    @('D24', '9.62'),
    @('E24', '  +3.52%  '),
    @('D25', '2.80'),
    @('E25', '  -5.07%  '),
    @('D27', '11.16'),
    @('E27', '  -1.68%  '),
    @('D28', '3.85'),
    @('E28', '  -1.98%  '),
    @('D29', '40.72'),
    @('E29', '  -0.29%  '),
    @('D30', '3.37'),
    @('E30', '  -2.28%  '),
    @('E31', '  +1.09%  '),
    @('D32', '170.82'),
    @('E32', '  -4.09%  '),
    @('D33', '21.08'),
    @('E33', '  -3.82%  '),
    @('D34', '0.0888'),
    @('E34', '  -2.51%  '),
    @('D35', '5.52'),
    @('E35', '  -1.39%  '),
    @('D36', '0.125'),
    @('E36', '  -1.48%  '),
    @('D37', '4.58'),
    @('E37', '  -5.51%  '),
    @('D38', '0.0345'),
    @('E38', '  -3.73%  '),
    @('E39', '  -9.18%  '),
    @('D40', '3.70'),
    @('E40', '  +0.36%  '),
    @('D41', '13.79'),
    @('E41', '  +12.98%  '),
    @('D42', '73.32'),
    @('E42', '  +9.33%  '),
    @('D43', '2.37'),
    @('E43', '  +1.77%  '),
    @('D44', '0.231'),
    @('E44', '  -2.48%  '),
    @('D45', '6.06'),
    @('E45', '  +10.64%  '),
    @('E46', '  +0.26%  '),
    @('D47', '1.35'),
    @('E47', '  -1.53%  '),
    @('D48', '8.48'),
    @('E48', '  -3.58%  '),
    @('D49', '0.0982'),
    @('E49', '  -3.65%  '),
    @('D50', '1.22'),
    @('E50', '  -0.52%  '),
    @('D51', '99.04'),
    @('E51', '  -0.44%  ')
)

foreach ($change in $changes) {
    $cellRef = $change[0]
    $newValue = $change[1]
    $trimmed = $newValue.Trim()
    $looksNumeric = $trimmed -match '^[+-]?\d+(\.\d+)?$'
    $range = $ws.Range($cellRef)
    if ($looksNumeric) {
        # Force text storage so e.g. '2.80' doesn't become the number 2.8
        $range.Value = "'" + $newValue
        $range.Style = $ws.Range('D26').Style
    } else {
        $range.Value = $newValue
    }
}
